$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.164.15'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.72%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.680.76'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '215.25'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.518'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +1.76%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '21.52'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +5.79%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0623'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.55%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0889'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.25%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.918.53'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.688.11'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.62%  '
$ws.Range('E14').Value = '  +1.54%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.538'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.81%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '66.35'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.75%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '27.153.11'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.70%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '238.59'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.61%  '
$ws.Range('E19').Value = '  +0.76%  '
$ws.Range('E20').Value = '  +1.32%  '
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.54'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +2.12%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.46'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.77%  '
$ws.Range('E24').Value = '  -3.72%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '148.26'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.98%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.26'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.34'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.77%  '
$ws.Range('E28').Value = '  +0.55%  '
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('E30').Value = '  +0.25%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.18'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.33%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.571.77'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +5.55%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.37'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.37%  '
$ws.Range('E34').Value = '  +2.34%  '
$ws.Range('E35').Value = '  +0.98%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.603'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +2.33%  '
$ws.Range('E37').Value = '  -0.98%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.935'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +4.07%  '
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('E40').Value = '  +3.63%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '69.19'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +2.66%  '
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.63'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -4.21%  '
$ws.Range('E44').Value = '  -2.47%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.826.95'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.67%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.787'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.12%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '90.90'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.41%  '
$ws.Range('E48').Value = '  +3.24%  '
$ws.Range('E49').Value = '  +1.91%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.13'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +6.12%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.104'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.42%  '
